$d = $word.ActiveDocument

# Center-align the title paragraph and the author paragraph.
$d.Paragraphs.Item(1).Alignment = 1
$d.Paragraphs.Item(2).Alignment = 1

# Move the "_GoBack" bookmark from the end of "Liam Whorriskey" to just
# after "Sp" in "Sprint 5 Retrospective" (splitting that run in two).
$d.Bookmarks.Item("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range(2, 2))
